$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (2021年) to the table, mirroring the formatting of
# the preceding year rows (bold header-style cell in column A, thin border,
# centered/top-aligned).
$row = 13

$ws.Cells.Item($row, 1).Value = "2021年"
$ws.Cells.Item($row, 1).Font.Bold = $true
$ws.Cells.Item($row, 1).HorizontalAlignment = -4108
$ws.Cells.Item($row, 1).VerticalAlignment = -4160
$ws.Cells.Item($row, 1).Borders.LineStyle = 1

$ws.Cells.Item($row, 2).Value = 3
$ws.Cells.Item($row, 3).Value = 49
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 1
$ws.Cells.Item($row, 6).Value = 3
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = 23
$ws.Cells.Item($row, 9).Value = 14
$ws.Cells.Item($row, 12).Value = 22
$ws.Cells.Item($row, 13).Value = 2
$ws.Cells.Item($row, 14).Value = 4
$ws.Cells.Item($row, 15).Value = 124
